$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: "The Control & Computing System hardware consists of:"
#    -> set font size to 12pt (sz/szCs = 24 half-points) on the paragraph
#       mark and on every run.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.Font.Size = 12
$r1.Font.SizeBi = 12

# ---------------------------------------------------------------------------
# 2) Bold the three one-line section headers that are list items
#    ("Flight Control System", "Flight Computing System", "Telemetry").
# ---------------------------------------------------------------------------
function Set-ParaBold($para) {
    $rng = $para.Range
    $rng.Font.Bold = 1
    $rng.Font.BoldBi = 1
}

$d.Paragraphs(2).Range.Find.Execute("Flight Control System") | Out-Null
Set-ParaBold $d.Paragraphs(2)
Set-ParaBold $d.Paragraphs(5)
Set-ParaBold $d.Paragraphs(8)

# ---------------------------------------------------------------------------
# 3) Remove the stray empty paragraph right after the telemetry bullet list,
#    and give the following paragraph a 0.5in first-line indent.
# ---------------------------------------------------------------------------
$pEmpty = $d.Paragraphs(10)
$pEmpty.Range.Delete()

$pHardware = $d.Paragraphs(10)
$pHardware.Range.ParagraphFormat.FirstLineIndent = 36

# ---------------------------------------------------------------------------
# 4) Bold "Complete autonomy of UAV:" and indent the paragraph under it.
# ---------------------------------------------------------------------------
Set-ParaBold $d.Paragraphs(11)
$d.Paragraphs(12).Range.ParagraphFormat.FirstLineIndent = 36

# ---------------------------------------------------------------------------
# 5) Bold "Advanced Computing & Artificial Intelligence:" and indent the
#    paragraph under it.
# ---------------------------------------------------------------------------
Set-ParaBold $d.Paragraphs(14)
$d.Paragraphs(15).Range.ParagraphFormat.FirstLineIndent = 36

# ---------------------------------------------------------------------------
# 6) "Flight Computer Software:" -> "Flight Controller Software:", bolded,
#    split across three runs ("Flight Co" / "ntroller" / " Software:") to
#    mirror the tracked edit.
# ---------------------------------------------------------------------------
$pSoft = $d.Paragraphs(16)
$softStart = $pSoft.Range.Start

# Bold the whole paragraph first (captures the paragraph-mark formatting).
$pSoft.Range.Font.Bold = 1
$pSoft.Range.Font.BoldBi = 1

# "Flight Co" + "mputer" + " Software:"  ->  replace "mputer" with "ntroller"
$midRng = $d.Range($softStart + 9, $softStart + 15)
$midRng.Text = "ntroller"

$pSoft2 = $d.Paragraphs(16)

$seg1 = $pSoft2.Range
$seg1.End = $softStart + 9
$seg1.Font.Bold = 0
$seg1.Font.Bold = 1
$seg1.Font.BoldBi = 0
$seg1.Font.BoldBi = 1

$seg2 = $pSoft2.Range
$seg2.Start = $softStart + 9
$seg2.End = $softStart + 17
$seg2.Font.Bold = 0
$seg2.Font.Bold = 1
$seg2.Font.BoldBi = 0
$seg2.Font.BoldBi = 1

$seg3 = $pSoft2.Range
$seg3.Start = $softStart + 17
$seg3.Font.Bold = 0
$seg3.Font.Bold = 1
$seg3.Font.BoldBi = 0
$seg3.Font.BoldBi = 1
